# Fruta / hortaliza, semanal
#
# This sheet ("Hortaliza, Terminal La Palmera de La Serena - Camote") holds a
# rolling weekly extract. On each weekly refresh the row data (date, volume,
# min/max/avg prices, $/Kg) gets re-pulled, which reshuffles which record
# lands on which row while columns A, B, C, E, F, G, H, I, N, O, Q, R stay
# fixed (same market/category/quality/packaging for every row). Only columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) change per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: RowNumber, D(Fecha), J(Volumen), K(PrecioMin), L(PrecioMax), M(PrecioProm), P(Precio$/Kg)
$data = @(
    @(2, 45154, 500, 16500, 17000, 16750, 931),
    @(3, 45166, 200, 16000, 17000, 16500, 917),
    @(4, 44957, 400, 21000, 22000, 21500, 1194),
    @(5, 44557, 400, 13000, 14000, 13500, 750),
    @(6, 45159, 400, 16000, 17000, 16500, 917),
    @(7, 45177, 540, 16000, 17000, 16500, 917),
    @(8, 44998, 320, 17000, 18000, 17500, 972),
    @(9, 44977, 400, 16500, 17000, 16750, 931),
    @(10, 45117, 300, 17000, 18000, 17500, 972),
    @(11, 44984, 200, 17000, 18000, 17500, 972),
    @(12, 44964, 300, 20000, 21000, 20500, 1139),
    @(13, 45180, 400, 16500, 17000, 16750, 931),
    @(14, 44568, 500, 15000, 16000, 15500, 861),
    @(15, 45142, 400, 17000, 18000, 17500, 972),
    @(16, 44547, 200, 13000, 14000, 13500, 750),
    @(17, 45068, 400, 16000, 17000, 16500, 917),
    @(18, 45005, 200, 17000, 18000, 17500, 972),
    @(19, 45194, 400, 16500, 17000, 16750, 931),
    @(20, 45152, 500, 16000, 17000, 16500, 917),
    @(21, 44960, 400, 19500, 20000, 19750, 1097)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]   # P - Precio $/Kg
}
